# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to match the freshly repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -1
    5  = -2
    6  = 9
    7  = 7
    8  = 4
    9  = 2
    10 = 2
    11 = -6
    12 = 1
    13 = -2
    14 = -1
    15 = 5
    17 = -1
    18 = 6
    21 = 10
    22 = 2
    23 = -3
    24 = 1
    25 = -3
    26 = 0
    27 = 2
    28 = -4
    29 = -4
    31 = -3
    32 = -1
    33 = 4
    35 = -1
    36 = -3
    37 = -4
    38 = 1
    39 = -1
    40 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
